# Correlation-Based Multivariate Stress Testing update
# - Append new trade rows (ALAB, AMD, AMZN) to the trade log on sheet 1 (工作表1)
# - Correct the traded-total-value figure for the existing row 173 trade
# - Update the window scroll/selection state to match the authored workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix existing D173 traded value (48647.05 -> 48694.400000000001) ---
$ws.Range("D173").Value = 48694.400000000001

# --- New trade rows appended after row 175 ---
# Carry over the formatting (date format on col A, number format on col D)
# from the row above by copy/paste-special so no new cell styles are created.
$ws.Range("A175").Copy()
$ws.Range("A176:A178").PasteSpecial(-4122)
$ws.Range("D175").Copy()
$ws.Range("D176:D178").PasteSpecial(-4122)

# Row 176: ALAB US Equity
$ws.Range("A176").Value = 46044
$ws.Range("B176").Value = "ALAB US Equity"
$ws.Range("C176").Value = -400
$ws.Range("D176").Value = 70534.759999999995

# Row 177: AMD US Equity
$ws.Range("A177").Value = 46044
$ws.Range("B177").Value = "AMD US Equity"
$ws.Range("C177").Value = 200
$ws.Range("D177").Value = -50490.54

# Row 178: AMZN US Equity
$ws.Range("A178").Value = 46044
$ws.Range("B178").Value = "AMZN US Equity"
$ws.Range("C178").Value = 300
$ws.Range("D178").Value = -70207.02

# --- View state: scroll position + active selection ---
$excel.ActiveWindow.ScrollRow = 156
$ws.Range("D173").Select()

$wb.Save()
